$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @(
    4145.857943051519,
    4071.457736692125,
    4024.627038190148,
    4007.831702098206,
    4007.831702098206,
    4007.831702098206,
    4007.831702098206,
    4007.831702098206,
    4007.831702098206,
    4007.831702098206,
    4007.831702098206
)

$row = 2
foreach ($val in $values) {
    $ws.Cells.Item($row, 3).Value = $val
    $row++
}
